$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3518.1
$ws.Range("I28").Value = 1901.375
$ws.Range("K28").Value = 1901.375
$ws.Range("M28").Value = -1416.375

$ws.Range("H116").Value = 11766966
$ws.Range("I116").Value = 33334784
$ws.Range("J116").Value = 2702.3635
$ws.Range("K116").Value = 33334784
$ws.Range("L116").Value = 2702.3635
$ws.Range("M116").Value = -33331342
$ws.Range("N116").Value = -9586.363499999999

$ws.Range("H125").Value = 3537
$ws.Range("I125").Value = 3954.3333
$ws.Range("K125").Value = 35588.9997
$ws.Range("M125").Value = -33128.9997

$ws.Range("H138").Value = 3740.5642
$ws.Range("I138").Value = 2148.923
$ws.Range("K138").Value = 6446.768999999999
$ws.Range("M138").Value = -1306.768999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 33022.6
$ws.Range("I10").Value = 28828.25
$ws.Range("J10").Value = 49800
$ws.Range("K10").Value = 28828.25
$ws.Range("L10").Value = 49800
$ws.Range("M10").Value = -28658.25
$ws.Range("N10").Value = -50140

$ws.Range("H45").Value = 1080.2222
$ws.Range("I45").Value = 1080.2222
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1080.2222
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -703.2221999999999
$ws.Range("N45").ClearContents() | Out-Null

$ws.Range("H122").Value = 2355.2144
$ws.Range("I122").Value = 2355.2144
$ws.Range("K122").Value = 7065.6432
$ws.Range("M122").Value = -4615.6432

$ws.Range("H135").Value = 35038.668
$ws.Range("J135").Value = 35038.668
$ws.Range("L135").Value = 35038.668
$ws.Range("N135").Value = -45178.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 11474.75
$ws.Range("I5").Value = 5999.7144
$ws.Range("J5").Value = 49800
$ws.Range("K5").Value = 5999.7144
$ws.Range("L5").Value = 49800
$ws.Range("M5").Value = -5886.7144
$ws.Range("N5").Value = -50026

$ws.Range("H62").Value = 29333.334
$ws.Range("J62").Value = 29333.334
$ws.Range("L62").Value = 29333.334
$ws.Range("N62").Value = -30705.334

$ws.Range("H65").Value = 29333.334
$ws.Range("J65").Value = 29333.334
$ws.Range("L65").Value = 88000.00199999999
$ws.Range("N65").Value = -94864.00199999999

$ws.Range("H134").Value = 3530.0833
$ws.Range("I134").Value = 3116.4666
$ws.Range("J134").Value = 4219.4443
$ws.Range("K134").Value = 9349.399800000001
$ws.Range("L134").Value = 12658.3329
$ws.Range("M134").Value = -6814.399800000001
$ws.Range("N134").Value = -17728.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 62843
$ws.Range("J52").Value = 62843
$ws.Range("L52").Value = 62843
$ws.Range("N52").Value = -63431

$ws.Range("H58").Value = 3365
$ws.Range("I58").Value = 3417
$ws.Range("J58").Value = 3295.6667
$ws.Range("K58").Value = 3417
$ws.Range("L58").Value = 3295.6667
$ws.Range("M58").Value = -3214
$ws.Range("N58").Value = -3701.6667

$ws.Range("H62").Value = 52761
$ws.Range("I62").Value = 64988.75
$ws.Range("J62").Value = 3850
$ws.Range("K62").Value = 64988.75
$ws.Range("L62").Value = 3850
$ws.Range("M62").Value = -64364.75
$ws.Range("N62").Value = -5098

$ws.Range("H65").Value = 52761
$ws.Range("I65").Value = 64988.75
$ws.Range("J65").Value = 3850
$ws.Range("K65").Value = 324943.75
$ws.Range("L65").Value = 19250
$ws.Range("M65").Value = -321823.75
$ws.Range("N65").Value = -25490

$ws.Range("H132").Value = 2441.2693
$ws.Range("I132").Value = 1938.7
$ws.Range("J132").Value = 4116.5
$ws.Range("K132").Value = 5816.1
$ws.Range("L132").Value = 12349.5
$ws.Range("M132").Value = -3286.1
$ws.Range("N132").Value = -17409.5

$ws.Range("H134").Value = 2012.8823
$ws.Range("I134").Value = 1579.9131
$ws.Range("J134").Value = 2918.182
$ws.Range("K134").Value = 4739.7393
$ws.Range("L134").Value = 8754.545999999998
$ws.Range("M134").Value = -2204.7393
$ws.Range("N134").Value = -13824.546

$ws.Range("H136").Value = 3365
$ws.Range("I136").Value = 3417
$ws.Range("J136").Value = 3295.6667
$ws.Range("K136").Value = 10251
$ws.Range("L136").Value = 9887.000100000001
$ws.Range("M136").Value = -7701
$ws.Range("N136").Value = -14987.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 533.7143
$ws.Range("I2").Value = 1023.9
$ws.Range("J2").Value = 88.09090999999999
$ws.Range("K2").Value = 6143.4
$ws.Range("L2").Value = 528.5454599999999
$ws.Range("M2").Value = -6030.4
$ws.Range("N2").Value = -754.5454599999999

$ws.Range("H15").Value = 259
$ws.Range("I15").Value = 51
$ws.Range("J15").Value = 363
$ws.Range("K15").Value = 153
$ws.Range("L15").Value = 1089
$ws.Range("M15").Value = -13
$ws.Range("N15").Value = -1369

$ws.Range("H113").Value = 628.70966
$ws.Range("I113").Value = 599.48
$ws.Range("J113").Value = 648.4595
$ws.Range("K113").Value = 1798.44
$ws.Range("L113").Value = 1945.3785
$ws.Range("M113").Value = 371.5599999999999
$ws.Range("N113").Value = -6285.378500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 9325.556
$ws.Range("J123").Value = 9325.556
$ws.Range("L123").Value = 9325.556
$ws.Range("N123").Value = -14225.556

$ws.Range("H132").Value = 3773.2122
$ws.Range("I132").Value = 2802
$ws.Range("J132").Value = 4137.4165
$ws.Range("K132").Value = 8406
$ws.Range("L132").Value = 12412.2495
$ws.Range("M132").Value = -5876
$ws.Range("N132").Value = -17472.2495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3140.6667
$ws.Range("I136").Value = 2257.4
$ws.Range("J136").Value = 3943.6365
$ws.Range("K136").Value = 6772.200000000001
$ws.Range("L136").Value = 11830.9095
$ws.Range("M136").Value = -4222.200000000001
$ws.Range("N136").Value = -16930.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 48999
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 48999
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 48999
$ws.Range("M8").ClearContents() | Out-Null
$ws.Range("N8").Value = -49279

$ws.Range("H81").Value = 90106.84
$ws.Range("I81").Value = 161827
$ws.Range("J81").Value = 6433.3335
$ws.Range("K81").Value = 323654
$ws.Range("L81").Value = 12866.667
$ws.Range("M81").Value = -322593
$ws.Range("N81").Value = -14988.667

$ws.Range("H84").Value = 90106.84
$ws.Range("I84").Value = 161827
$ws.Range("J84").Value = 6433.3335
$ws.Range("K84").Value = 1618270
$ws.Range("L84").Value = 64333.335
$ws.Range("M84").Value = -1612966
$ws.Range("N84").Value = -74941.33499999999

$ws.Range("H122").Value = 43105950
$ws.Range("I122").Value = 69446020
$ws.Range("K122").Value = 208338060
$ws.Range("M122").Value = -208335610

$ws.Range("H132").Value = 5222.6665
$ws.Range("I132").Value = 5401
$ws.Range("J132").Value = 4999.75
$ws.Range("K132").Value = 16203
$ws.Range("L132").Value = 14999.25
$ws.Range("M132").Value = -13673
$ws.Range("N132").Value = -20059.25

$ws.Range("H136").Value = 2457.1785
$ws.Range("I136").Value = 2135.1
$ws.Range("J136").Value = 3262.375
$ws.Range("K136").Value = 6405.299999999999
$ws.Range("L136").Value = 9787.125
$ws.Range("M136").Value = -3855.299999999999
$ws.Range("N136").Value = -14987.125
